$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rerun DLinear on new dataset: replace the validation-scaler close series
# (index in col A, scaled close value in col B). The new run produced 132
# rows (index 0..131) instead of the previous 126 (index 0..125), so the
# sheet grows from A1:B127 to A1:B133.
$data = @(
    @(0, 1.19685),
    @(1, 1.16442),
    @(2, 1.15956),
    @(3, 1.17902),
    @(4, 1.19685),
    @(5, 1.16118),
    @(6, 1.16767),
    @(7, 1.15956),
    @(8, 1.16605),
    @(9, 1.17577),
    @(10, 1.17902),
    @(11, 1.18713),
    @(12, 1.20658),
    @(13, 1.22604),
    @(14, 1.22928),
    @(15, 1.21793),
    @(16, 1.21631),
    @(17, 1.22928),
    @(18, 1.25847),
    @(19, 1.23253),
    @(20, 1.24874),
    @(21, 1.24874),
    @(22, 1.21307),
    @(23, 1.2228),
    @(24, 1.22442),
    @(25, 1.23091),
    @(26, 1.2228),
    @(27, 1.2455),
    @(28, 1.26496),
    @(29, 1.26009),
    @(30, 1.25523),
    @(31, 1.24712),
    @(32, 1.24388),
    @(33, 1.25036),
    @(34, 1.26496),
    @(35, 1.2682),
    @(36, 1.26658),
    @(37, 1.25523),
    @(38, 1.22442),
    @(39, 1.22928),
    @(40, 1.22442),
    @(41, 1.23577),
    @(42, 1.24064),
    @(43, 1.24064),
    @(44, 1.23901),
    @(45, 1.25036),
    @(46, 1.25199),
    @(47, 1.23901),
    @(48, 1.22928),
    @(49, 1.23091),
    @(50, 1.24388),
    @(51, 1.24388),
    @(52, 1.23901),
    @(53, 1.24064),
    @(54, 1.30874),
    @(55, 1.34279),
    @(56, 1.36712),
    @(57, 1.37684),
    @(58, 1.38982),
    @(59, 1.37522),
    @(60, 1.39306),
    @(61, 1.39144),
    @(62, 1.38009),
    @(63, 1.40603),
    @(64, 1.38495),
    @(65, 1.3736),
    @(66, 1.44008),
    @(67, 1.43846),
    @(68, 1.44819),
    @(69, 1.4563),
    @(70, 1.48549),
    @(71, 1.49522),
    @(72, 1.48387),
    @(73, 1.45792),
    @(74, 1.49359),
    @(75, 1.47251),
    @(76, 1.47251),
    @(77, 1.50332),
    @(78, 1.58602),
    @(79, 1.58602),
    @(80, 1.56981),
    @(81, 1.52116),
    @(82, 1.50495),
    @(83, 1.54548),
    @(84, 1.57305),
    @(85, 1.55197),
    @(86, 1.53738),
    @(87, 1.55359),
    @(88, 1.55683),
    @(89, 1.58764),
    @(90, 1.57791),
    @(91, 1.58278),
    @(92, 1.58602),
    @(93, 1.55521),
    @(94, 1.53738),
    @(95, 1.53251),
    @(96, 1.51467),
    @(97, 1.52116),
    @(98, 1.53413),
    @(99, 1.55197),
    @(100, 1.5617),
    @(101, 1.49846),
    @(102, 1.52116),
    @(103, 1.5017),
    @(104, 1.4563),
    @(105, 1.47738),
    @(106, 1.50981),
    @(107, 1.63629),
    @(108, 1.68656),
    @(109, 1.68656),
    @(110, 1.75304),
    @(111, 1.73034),
    @(112, 1.7579),
    @(113, 1.8179),
    @(114, 1.81304),
    @(115, 1.80169),
    @(116, 1.81304),
    @(117, 1.79194),
    @(118, 1.80493),
    @(119, 1.8633),
    @(120, 1.88438),
    @(121, 1.86979),
    @(122, 1.84547),
    @(123, 1.89411),
    @(124, 1.92654),
    @(125, 1.91844),
    @(126, 1.82763),
    @(127, 1.84222),
    @(128, 1.91033),
    @(129, 1.88925),
    @(130, 1.8779),
    @(131, 1.87141)
)

$startRow = 2
$oldLastRow = 127
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
}

# New rows beyond the old extent need the index-column formatting (bold,
# bordered, centered) that row 127 already carries, since brand-new cells
# otherwise come back with the default style.
$endRow = $startRow + $data.Length - 1
if ($endRow -gt $oldLastRow) {
    $ws.Range("A$oldLastRow").Copy()
    $ws.Range("A" + ($oldLastRow + 1) + ":A$endRow").PasteSpecial(-4122)
}
